# develop-master : fix import karyawan
#
# The template's header row (and the 5 sample data rows beneath it) is
# missing several columns that the importer now expects. Insert them in
# their correct positions, relative to the existing "kolom" layout:
#
#   ... nik | [kd_entitas] | ket_jabatan | kd_subdiv | [kd_bagian] | id_cabang | kd_jabatan |
#   kd_pangkat_golongan | [status_jabatan] | id_is | ...
#   ... alamat_sekarang | [pendidikan] | [tgl_mulai] | [tanggal_penonaktifan] |
#   [sk_pemberhentian] | kpj | ...
#
# Columns are inserted from right to left so that earlier (left-hand)
# insertion points keep their original A1 addresses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 4 new columns before the old column R (kpj) -------------------------
$ws.Range("R1:U1").EntireColumn.Insert()
$ws.Range("R1").Value = "pendidikan"
$ws.Range("S1").Value = "tgl_mulai"
$ws.Range("T1").Value = "tanggal_penonaktifan"
$ws.Range("U1").Value = "sk_pemberhentian"
$ws.Columns.Item(18).ColumnWidth = 14.166666666666666
$ws.Columns.Item(19).ColumnWidth = 14.166666666666666
$ws.Columns.Item(20).ColumnWidth = 14.166666666666666
$ws.Columns.Item(21).ColumnWidth = 14.166666666666666

# --- 1 new column before the old column I (id_is) -------------------------
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("I1").Value = "status_jabatan"
$ws.Columns.Item(9).ColumnWidth = 18.166666666666668

# --- 1 new column before the old column F (kd_pangkat_golongan) -----------
$ws.Range("F1").EntireColumn.Insert()
$ws.Range("F1").Value = "kd_bagian"
$ws.Columns.Item(6).ColumnWidth = 9.5

# --- 1 new column before the old column D (ket_jabatan) -------------------
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("D1").Value = "kd_entitas"
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666

# --- put the selection/view back roughly where the authored file had it ---
$ws.Range("X1").Select()
